$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(6, 11).Value = 45396
$ws.Cells.Item(7, 11).Value = 45396
$ws.Cells.Item(8, 11).Value = 45396
$ws.Cells.Item(9, 11).Value = 45396
$ws.Cells.Item(10, 11).Value = 45396
$ws.Cells.Item(11, 11).Value = 45396
$ws.Cells.Item(12, 11).Value = 45398
$ws.Cells.Item(14, 11).Value = 45399
$ws.Cells.Item(15, 11).Value = 45399
$ws.Cells.Item(16, 11).Value = 45399
$ws.Cells.Item(17, 11).Value = 45399
$ws.Cells.Item(18, 11).Value = 45399
$ws.Cells.Item(19, 11).Value = 45399
$ws.Cells.Item(20, 11).Value = 45399
$ws.Cells.Item(21, 11).Value = 45400
$ws.Cells.Item(22, 11).Value = 45448
$ws.Cells.Item(24, 11).Value = 45402
$ws.Cells.Item(25, 11).Value = 45402
$ws.Cells.Item(26, 11).Value = 45402
$ws.Cells.Item(27, 11).Value = 45402
$ws.Cells.Item(28, 11).Value = 45403
$ws.Cells.Item(29, 11).Value = 45403
$ws.Cells.Item(30, 11).Value = 45403
$ws.Cells.Item(31, 11).Value = 45403
$ws.Cells.Item(32, 11).Value = 45403
$ws.Cells.Item(34, 11).Value = 45409
$ws.Cells.Item(35, 11).Value = 45409
$ws.Cells.Item(36, 11).Value = 45409
$ws.Cells.Item(37, 11).Value = 45410
$ws.Cells.Item(38, 11).Value = 45410
$ws.Cells.Item(39, 11).Value = 45410
$ws.Cells.Item(40, 11).Value = 45410
$ws.Cells.Item(41, 11).Value = 45410
$ws.Cells.Item(42, 11).Value = 45410
$ws.Cells.Item(43, 11).Value = 45411
$ws.Cells.Item(44, 11).Value = 45416
$ws.Cells.Item(45, 11).Value = 45416
$ws.Cells.Item(46, 11).Value = 45416
$ws.Cells.Item(47, 11).Value = 45417
$ws.Cells.Item(48, 11).Value = 45417
$ws.Cells.Item(49, 11).Value = 45417
$ws.Cells.Item(50, 11).Value = 45417
$ws.Cells.Item(51, 11).Value = 45448
$ws.Cells.Item(52, 11).Value = 45423
$ws.Cells.Item(53, 11).Value = 45424
$ws.Cells.Item(54, 11).Value = 45424
$ws.Cells.Item(55, 11).Value = 45424
$ws.Cells.Item(56, 11).Value = 45424
$ws.Cells.Item(57, 11).Value = 45424
$ws.Cells.Item(58, 11).Value = 45425
$ws.Cells.Item(59, 11).Value = 45452
$ws.Cells.Item(60, 11).Value = 45444
$ws.Cells.Item(61, 11).Value = 45444
$ws.Cells.Item(62, 11).Value = 45444
$ws.Cells.Item(63, 11).Value = 45444
$ws.Cells.Item(64, 11).Value = 45444
$ws.Cells.Item(65, 11).Value = 45445
$ws.Cells.Item(66, 11).Value = 45445
$ws.Cells.Item(67, 11).Value = 45445
$ws.Cells.Item(68, 11).Value = 45445
$ws.Cells.Item(69, 11).Value = 45445
$ws.Cells.Item(70, 11).Value = 45454
$ws.Cells.Item(71, 11).Value = 45454
$ws.Cells.Item(72, 11).Value = 45454
$ws.Cells.Item(73, 11).Value = 45454
$ws.Cells.Item(74, 11).Value = 45456
$ws.Cells.Item(75, 11).Value = 45456
$ws.Cells.Item(76, 11).Value = 45456
$ws.Cells.Item(77, 11).Value = 45456
$ws.Cells.Item(78, 11).Value = 45456
$ws.Cells.Item(79, 11).Value = 45456
$ws.Cells.Item(80, 11).Value = 45458
$ws.Cells.Item(81, 11).Value = 45458
$ws.Cells.Item(82, 11).Value = 45459
$ws.Cells.Item(83, 11).Value = 45459
$ws.Cells.Item(84, 11).Value = 45459
$ws.Cells.Item(85, 11).Value = 45459
$ws.Cells.Item(86, 11).Value = 45459
$ws.Cells.Item(87, 11).Value = 45459
$ws.Cells.Item(88, 11).Value = 45459
$ws.Cells.Item(89, 11).Value = 45460
$ws.Cells.Item(90, 11).Value = 45462
$ws.Cells.Item(91, 11).Value = 45462
$ws.Cells.Item(92, 11).Value = 45462
$ws.Cells.Item(93, 11).Value = 45462
$ws.Cells.Item(94, 11).Value = 45462
$ws.Cells.Item(95, 11).Value = 45462
$ws.Cells.Item(96, 11).Value = 45462
$ws.Cells.Item(97, 11).Value = 45463
$ws.Cells.Item(98, 11).Value = 45463
$ws.Cells.Item(99, 11).Value = 45463
$ws.Cells.Item(100, 11).Value = 45465
$ws.Cells.Item(101, 11).Value = 45465
$ws.Cells.Item(102, 11).Value = 45465
$ws.Cells.Item(103, 11).Value = 45465
$ws.Cells.Item(104, 11).Value = 45466
$ws.Cells.Item(105, 11).Value = 45466
$ws.Cells.Item(106, 11).Value = 45466
$ws.Cells.Item(107, 11).Value = 45466
$ws.Cells.Item(108, 11).Value = 45466
$ws.Cells.Item(109, 11).Value = 45466
$ws.Cells.Item(110, 11).Value = 45469
$ws.Cells.Item(111, 11).Value = 45469
$ws.Cells.Item(112, 11).Value = 45469
$ws.Cells.Item(113, 11).Value = 45469
$ws.Cells.Item(114, 11).Value = 45469
$ws.Cells.Item(115, 11).Value = 45469
$ws.Cells.Item(116, 11).Value = 45469
$ws.Cells.Item(117, 11).Value = 45469
$ws.Cells.Item(118, 11).Value = 45470
$ws.Cells.Item(119, 11).Value = 45470
$ws.Cells.Item(120, 11).Value = 45472
$ws.Cells.Item(121, 11).Value = 45472
$ws.Cells.Item(122, 11).Value = 45473
$ws.Cells.Item(123, 11).Value = 45473
$ws.Cells.Item(124, 11).Value = 45473
$ws.Cells.Item(125, 11).Value = 45473
$ws.Cells.Item(126, 11).Value = 45473
$ws.Cells.Item(127, 11).Value = 45473
$ws.Cells.Item(128, 11).Value = 45473
$ws.Cells.Item(129, 11).Value = 45474
$ws.Cells.Item(130, 11).Value = 45476
$ws.Cells.Item(131, 11).Value = 45476
$ws.Cells.Item(132, 11).Value = 45476
$ws.Cells.Item(133, 11).Value = 45476
$ws.Cells.Item(134, 11).Value = 45476
$ws.Cells.Item(135, 11).Value = 45476
$ws.Cells.Item(136, 11).Value = 45477
$ws.Cells.Item(137, 11).Value = 45477
$ws.Cells.Item(138, 11).Value = 45477
$ws.Cells.Item(139, 11).Value = 45477
$ws.Cells.Item(140, 11).Value = 45479
$ws.Cells.Item(141, 11).Value = 45479
$ws.Cells.Item(142, 11).Value = 45480
$ws.Cells.Item(143, 11).Value = 45480
$ws.Cells.Item(144, 11).Value = 45480
$ws.Cells.Item(145, 11).Value = 45480
$ws.Cells.Item(146, 11).Value = 45480
$ws.Cells.Item(147, 11).Value = 45480
$ws.Cells.Item(148, 11).Value = 45480
$ws.Cells.Item(149, 11).Value = 45480
$ws.Cells.Item(151, 11).Value = 45483
$ws.Cells.Item(152, 11).Value = 45483
$ws.Cells.Item(154, 11).Value = 45483
$ws.Cells.Item(155, 11).Value = 45484
$ws.Cells.Item(156, 11).Value = 45484
$ws.Cells.Item(157, 11).Value = 45484
$ws.Cells.Item(158, 11).Value = 45484
$ws.Cells.Item(159, 11).Value = 45484
